# Updated the results spreadsheet with an accurate Kakuro experiment
#
# Row 9 (the Kakuro "conceptis" entry) previously recorded only a Total of 1
# with no Hits value. The real experiment numbers are Hits=7, Total=11, and
# a new Notes-style remark is added in column H explaining why a few of the
# steps weren't picked up by the solver.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 11
$ws.Range("H9").Value = "step 7 is done by us way earlier, step 8 is a MUS of size 8, step 9 and 10 we find simpler/other deductions first that leave these as super easy deductions."

# Move the active selection to match the author's saved cursor position.
[void]$ws.Range("C9").Select()
